# "Generate Report for Archive"
#  - Update status text "Ready for handoff" -> "In Translation" on all sheets
#  - Shrink the "Status"/"zh-cn"/"de-de" columns to their new auto-fit width

$wb = $excel.ActiveWorkbook

# 1) Replace the status text everywhere it appears (Overview, zh-cn, de-de)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# 2) Narrow the previously-widened columns to match the new auto-fit width
#    (raw OOXML width 17.2159881591797 -> 13.4101848602295, i.e.
#    ColumnWidth 12.5 which this engine reports/stores as ~13.3333)
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de status)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
